$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("terralith:cave/underground_jungle", 30),
    @("terralith:cave/thermal_caves", 30),
    @("terralith:cave/infested_caves", 10),
    @("terralith:cave/mantle_caves", 60),
    @("the_winter_rescue:magmatic_deposits", 50),
    @("the_winter_rescue:hydrothermal_deposits", 40)
)

$startRow = 32
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $newRows[$i][1]
}
